# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Values below were recalculated (regen std/mean, calc and write s_vals)
# and differ from the previously stored figures; only column G (header "K")
# needs updating, rows 2 through 36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(7, 5, 2, 6, 6, 9, 8, 0, 9, 6, 2, 6, 2, 5, 9, 2, 5, 6, 8, 9, 4, 7, 10, 7, 7, 1, 7, 3, 7, 1, 7, 3, 5, 5, 5)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
